# Append one row (row 10) of data to Sheet1, matching the source system's
# convention of storing every field (even numeric-looking ones) as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and H contain values that Excel could otherwise auto-convert
# (a digits-only string -> number, a slash/colon-laden string -> date/time).
# Force them to Text format first so the literal string is preserved.
$ws.Range("C10").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"

$ws.Range("A10").Value = ""
$ws.Range("B10").Value = "أحمد شريم"
$ws.Range("C10").Value = "2323"
$ws.Range("D10").Value = "ايتا"
$ws.Range("E10").Value = "الرحلة 2"
$ws.Range("F10").Value = "C2"
$ws.Range("G10").Value = "NRC"
$ws.Range("H10").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٢٤:٤٧ م"
